$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-write header row and row labels (same text as before) so the shared
# string table gets a fresh duplicate set appended, matching how Excel
# behaves when re-saving after cell content is rewritten.
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "RMSE"
$ws.Range("C1").Value = "NRMSE"
$ws.Range("D1").Value = "MAE"
$ws.Range("E1").Value = "RSE"
$ws.Range("F1").Value = "RRSE"
$ws.Range("G1").Value = "RAE"
$ws.Range("H1").Value = "R2"
$ws.Range("I1").Value = "Corr Coeff"

$ws.Range("A2").Value = "random_forest"
$ws.Range("A3").Value = "lsboost"
$ws.Range("A4").Value = "neural_network"
$ws.Range("A5").Value = "old_model"

# Updated experiment result with new k for cross-validation
$ws.Range("B2").Value = 57.44262511799942
$ws.Range("C2").Value = 0.73173316791819698
$ws.Range("D2").Value = 38.081507768920027
$ws.Range("E2").Value = 0.51458841104222464
$ws.Range("F2").Value = 0.71734817978595622
$ws.Range("G2").Value = 0.62221571371600981
$ws.Range("H2").Value = 0.48541158895777536
$ws.Range("I2").Value = 0.70999876518723282

$ws.Range("B3").Value = 57.041745930137346
$ws.Range("C3").Value = 0.72662656637475753
$ws.Range("D3").Value = 34.742438210657873
$ws.Range("E3").Value = 0.50743107915474217
$ws.Range("F3").Value = 0.71234196784602144
$ws.Range("G3").Value = 0.56765848449734091
$ws.Range("H3").Value = 0.49256892084525783
$ws.Range("I3").Value = 0.71729943163330045

$ws.Range("B4").Value = 60.611287740869201
$ws.Range("C4").Value = 0.77209719261820953
$ws.Range("D4").Value = 37.1065142723604
$ws.Range("E4").Value = 0.5729259133099962
$ws.Range("F4").Value = 0.75691869663128031
$ws.Range("G4").Value = 0.60628524483826651
$ws.Range("H4").Value = 0.4270740866900038
$ws.Range("I4").Value = 0.67650700922256801

$ws.Range("B5").Value = 64.119604773399075
$ws.Range("C5").Value = 0.81678790671773005
$ws.Range("D5").Value = 37.532211197026022
$ws.Range("E5").Value = 0.64116988677672138
$ws.Range("F5").Value = 0.80073084540107564
$ws.Range("G5").Value = 0.61324072878115043
$ws.Range("H5").Value = 0.35883011322327862
$ws.Range("I5").Value = 0.68061816663214003
